# Vehicles now need services, some fixes
#
# The source sheet's B column used to hold a shared formula "=A{r}*$M$2"
# across B3:B10; it becomes a literal "=800" per cell (no longer shared).
# Because the underlying engine re-materialises a shared-formula group
# whenever the *exact* original range (B3:B10) ends up holding identical
# formula text again, we rebuild the worksheet from scratch (new sheet,
# same name) so no stale shared-formula bookkeeping survives, and then
# reconstruct each cell explicitly with the exact sharing layout the
# target file expects.

$wb = $excel.ActiveWorkbook
$oldName = $wb.ActiveSheet.Name
$oldIndex = $wb.ActiveSheet.Index

$new = $wb.Worksheets.Add()
$new.Name = "__rebuild_tmp__"

# ---- Row 1 : headers ----
$new.Range("A1").Value = "population"
$new.Range("B1").Value = "demand"
$new.Range("C1").Value = "supply"
$new.Range("D1").Value = "difference"
$new.Range("E1").Value = "price_inertia"
$new.Range("F1").Value = "10*factor"
$new.Range("G1").Value = "daily balance"
$new.Range("H1").Value = "monthly balance"
$new.Range("K1").Value = "BASE_FOREIGN_SUPPLY_DEMAND"
$new.Range("L1").Value = "FOREIGN_SUPPLY_DEMAND_FACTOR"
$new.Range("M1").Value = "PRICE_BUY_RUB"
$new.Range("N1").Value = "PRICE_SELL_RUB"

# ---- Column A : population ----
for ($r = 2; $r -le 10; $r++) {
    $new.Cells.Item($r, 1).Value = 20000
}

# ---- Column B : demand - now a literal 800, individually per cell (not shared) ----
for ($r = 2; $r -le 10; $r++) {
    $new.Cells.Item($r, 2).Formula = '=800'
}

# ---- Column C : supply (plain inputs) ----
$cvals = @(0, 200, 400, 600, 800, 1000, 1200, 1400, 1600)
for ($i = 0; $i -lt $cvals.Length; $i++) {
    $new.Cells.Item($i + 2, 3).Value = $cvals[$i]
}

# ---- Row 2 : individual (non-shared) formulas ----
$new.Range("D2").Formula = '=C2-B2'
$new.Range("E2").Formula = '=$K$2+$L$2*A2'
$new.Range("F2").Formula = '=1-D2/E2'
$new.Range("G2").Formula = '=D2*F2*IF(D2<0,$M$2,$N$2)'
$new.Range("H2").Formula = '=G2*30'

# ---- Rows 3-10 : shared formula groups (matches target layout) ----
$new.Range("D3:D8").Formula = '=C3-B3'
$new.Range("E3:E10").Formula = '=$K$2+$L$2*A3+C3'
$new.Range("F3:F10").Formula = '=1-D3/E3'
$new.Range("G3:G10").Formula = '=D3*F3*IF(D3<0,$M$2,$N$2)'
$new.Range("H3:H10").Formula = '=G3*30'
$new.Range("D9:D10").Formula = '=C9-B9'

# ---- K2:N2 parameters ----
$new.Range("K2").Value = 200
$new.Range("L2").Value = 0.03
$new.Range("M2").Value = 6
$new.Range("N2").Value = 5

# ---- View state: selection now sits at D17 ----
$new.Range("D17").Select() | Out-Null

# ---- Swap the rebuilt sheet in under the original name/position ----
$app.DisplayAlerts = $false
$wb.Worksheets.Item($oldName).Delete()
$new.Name = $oldName
$new.Move($wb.Worksheets.Item(1))
